$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "Priority" counter column ---
$ws.Columns.Item(7).ColumnWidth = 14

$ws.Range("G4").Value = "Priority"

$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 2
$ws.Range("G12").Value = 3
$ws.Range("G13").Value = 3

$ws.Range("G16").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("G18").Value = 2
$ws.Range("G19").Value = 1
$ws.Range("G20").Value = 1

# --- ARTICLES section: GET /articles/get/{article} -> GET /articles/{article} ---
$ws.Range("C20").Value = "/articles/{article}"

$ws.Range("G23").Value = 1
$ws.Range("G24").Value = 11
$ws.Range("G25").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("G27").Value = 1

# --- New CATEGORIES section header (row 29), styled like the other section headers ---
$ws.Range("B29").Value = "CATEGORIES"
$ws.Range("B29").Font.Bold = $true
$ws.Range("B29").Font.Color = 7884319

# --- New CATEGORIES rows ---
$ws.Range("B30").Value = "GET"
$ws.Range("C30").Value = "/category/list"
$ws.Range("D30").Value = "QueryDTO"
$ws.Range("E30").Value = "null"
$ws.Range("G30").Value = 1

$ws.Range("B31").Value = "GET"
$ws.Range("C31").Value = "/category/{id}"
$ws.Range("D31").Value = "null"
$ws.Range("E31").Value = "number"
$ws.Range("G31").Value = 1

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- View state: selection moved to D38 ---
$ws.Range("D38").Select()
